# Generate Report for Handoff
# Updates the localization-status workbook to reflect a new handoff report:
#  - new source files (ff33dccf... and ffff51ef5982...) replace the previous
#    pair, now in "Ready for handoff" status with a fresh handoff timestamp
#  - the per-language sheets lose their "Latest Target File"/"Latest Handback
#    File" columns' data (no handback has happened yet for the new files)
#  - hyperlink display text is refreshed to match the new file names

$wb = $excel.ActiveWorkbook

function Set-HyperlinkDisplay($ws, $row, $col, $text) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Row -eq $row -and $h.Range.Column -eq $col) {
            $h.TextToDisplay = $text
        }
    }
}

function Remove-HyperlinksAt($ws, $row, $cols) {
    $changed = $true
    while ($changed) {
        $changed = $false
        foreach ($h in $ws.Hyperlinks) {
            if ($h.Range.Row -eq $row -and ($cols -contains $h.Range.Column)) {
                $h.Delete()
                $changed = $true
                break
            }
        }
    }
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "ff33dccf-443a-4cfa-bf0b-2cea559efc06.md"
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-52-20 12:52:09"

$wsOverview.Range("A3").Value = "ffff51ef5982-defa-4d2c-ba7f-d6cf46eed4f6.md"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-52-20 12:52:09"

Set-HyperlinkDisplay $wsOverview 2 1 "ff33dccf-443a-4cfa-bf0b-2cea559efc06.md"
Set-HyperlinkDisplay $wsOverview 3 1 "ffff51ef5982-defa-4d2c-ba7f-d6cf46eed4f6.md"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "ff33dccf-443a-4cfa-bf0b-2cea559efc06.md"
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("D2").Value = "ff33dccf-443a-4cfa-bf0b-2cea559efc06.a854402d478cbeb40495a7cb4dba6d7f7f6bb74a.zh-cn.xlf"
$wsZh.Range("E2").Value = "2016-03-20 12:52:06"
$wsZh.Range("H2").Value = "0001-01-01 00:00:00"
$wsZh.Range("I2").Value = "Include"

$wsZh.Range("A3").Value = "ffff51ef5982-defa-4d2c-ba7f-d6cf46eed4f6.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "ff33dccf-443a-4cfa-bf0b-2cea559efc06.a854402d478cbeb40495a7cb4dba6d7f7f6bb74a.zh-cn.xlf"
$wsZh.Range("E3").Value = "2016-03-20 12:52:06"
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"
$wsZh.Range("I3").Value = "Include"

Remove-HyperlinksAt $wsZh 2 @(6, 7)
Remove-HyperlinksAt $wsZh 3 @(6, 7)
$wsZh.Range("F2:G3").Clear()

Set-HyperlinkDisplay $wsZh 2 1 "ff33dccf-443a-4cfa-bf0b-2cea559efc06.md"
Set-HyperlinkDisplay $wsZh 2 4 "ff33dccf-443a-4cfa-bf0b-2cea559efc06.a854402d478cbeb40495a7cb4dba6d7f7f6bb74a.zh-cn.xlf"
Set-HyperlinkDisplay $wsZh 3 1 "ffff51ef5982-defa-4d2c-ba7f-d6cf46eed4f6.md"
Set-HyperlinkDisplay $wsZh 3 4 "ff33dccf-443a-4cfa-bf0b-2cea559efc06.a854402d478cbeb40495a7cb4dba6d7f7f6bb74a.zh-cn.xlf"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "ff33dccf-443a-4cfa-bf0b-2cea559efc06.md"
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("D2").Value = "ff33dccf-443a-4cfa-bf0b-2cea559efc06.a854402d478cbeb40495a7cb4dba6d7f7f6bb74a.de-de.xlf"
$wsDe.Range("E2").Value = "2016-03-20 12:52:09"
$wsDe.Range("H2").Value = "0001-01-01 00:00:00"
$wsDe.Range("I2").Value = "Include"

$wsDe.Range("A3").Value = "ffff51ef5982-defa-4d2c-ba7f-d6cf46eed4f6.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "ff33dccf-443a-4cfa-bf0b-2cea559efc06.a854402d478cbeb40495a7cb4dba6d7f7f6bb74a.de-de.xlf"
$wsDe.Range("E3").Value = "2016-03-20 12:52:09"
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDe.Range("I3").Value = "Include"

Remove-HyperlinksAt $wsDe 2 @(6, 7)
Remove-HyperlinksAt $wsDe 3 @(6, 7)
$wsDe.Range("F2:G3").Clear()

Set-HyperlinkDisplay $wsDe 2 1 "ff33dccf-443a-4cfa-bf0b-2cea559efc06.md"
Set-HyperlinkDisplay $wsDe 2 4 "ff33dccf-443a-4cfa-bf0b-2cea559efc06.a854402d478cbeb40495a7cb4dba6d7f7f6bb74a.de-de.xlf"
Set-HyperlinkDisplay $wsDe 3 1 "ffff51ef5982-defa-4d2c-ba7f-d6cf46eed4f6.md"
Set-HyperlinkDisplay $wsDe 3 4 "ff33dccf-443a-4cfa-bf0b-2cea559efc06.a854402d478cbeb40495a7cb4dba6d7f7f6bb74a.de-de.xlf"
